$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 216.85715
$ws.Range("I9").Value = 78.2
$ws.Range("J9").Value = 563.5
$ws.Range("K9").Value = 78.2
$ws.Range("L9").Value = 563.5
$ws.Range("M9").Value = 90.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 4400
$ws.Range("I20").Value = 4400
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 4400
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -4170

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1334.8
$ws.Range("I32").Value = 1037.5
$ws.Range("J32").Value = 1533
$ws.Range("K32").Value = 1037.5
$ws.Range("L32").Value = 1533
$ws.Range("M32").Value = -711.5
$ws.Range("N32").Value = -2185

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 4400
$ws.Range("I35").Value = 4400
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4400
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -4021

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7370
$ws.Range("I43").Value = 5499.5
$ws.Range("J43").Value = 11111
$ws.Range("K43").Value = 5499.5
$ws.Range("L43").Value = 11111
$ws.Range("M43").Value = -5430.5
$ws.Range("N43").Value = -11249

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4985.5713
$ws.Range("I70").Value = 8633.333000000001
$ws.Range("J70").Value = 2249.75
$ws.Range("K70").Value = 25899.999
$ws.Range("L70").Value = 6749.25
$ws.Range("M70").Value = -25629.999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 4985.5713
$ws.Range("I73").Value = 8633.333000000001
$ws.Range("J73").Value = 2249.75
$ws.Range("K73").Value = 25899.999
$ws.Range("L73").Value = 6749.25
$ws.Range("M73").Value = -24963.999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 9614.615
$ws.Range("I113").Value = 12782
$ws.Range("J113").Value = 4546.8
$ws.Range("K113").Value = 12782
$ws.Range("L113").Value = 4546.8
$ws.Range("M113").Value = -9528
$ws.Range("N113").Value = -11054.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3500743.5
$ws.Range("I125").Value = 5687080
$ws.Range("J125").Value = 2605
$ws.Range("K125").Value = 51183720
$ws.Range("L125").Value = 23445
$ws.Range("M125").Value = -51181260
$ws.Range("N125").Value = -28365

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2208.5205
$ws.Range("I138").Value = 1505.3226
$ws.Range("J138").Value = 2533.8806
$ws.Range("K138").Value = 4515.9678
$ws.Range("L138").Value = 7601.641799999999
$ws.Range("M138").Value = 624.0321999999996
$ws.Range("N138").Value = -17881.6418

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3961.0715
$ws.Range("I141").Value = 4415
$ws.Range("J141").Value = 2826.25
$ws.Range("K141").Value = 13245
$ws.Range("L141").Value = 8478.75
$ws.Range("M141").Value = -8065

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3331.2632
$ws.Range("I32").Value = 2111.9795
$ws.Range("J32").Value = 10799.375
$ws.Range("K32").Value = 2111.9795
$ws.Range("L32").Value = 10799.375
$ws.Range("M32").Value = -1824.9795

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2485.4827
$ws.Range("I45").Value = 2517.2727
$ws.Range("J45").Value = 2385.5715
$ws.Range("K45").Value = 2517.2727
$ws.Range("L45").Value = 2385.5715
$ws.Range("M45").Value = -2140.2727
$ws.Range("N45").Value = -3139.5715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 14832.549
$ws.Range("I74").Value = 1951.48
$ws.Range("J74").Value = 68503.664
$ws.Range("K74").Value = 1951.48
$ws.Range("L74").Value = 68503.664
$ws.Range("M74").Value = -1077.48
$ws.Range("N74").Value = -70251.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 14832.549
$ws.Range("I77").Value = 1951.48
$ws.Range("J77").Value = 68503.664
$ws.Range("K77").Value = 9757.4
$ws.Range("L77").Value = 342518.32
$ws.Range("M77").Value = -5389.4
$ws.Range("N77").Value = -351254.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1122881.9
$ws.Range("I122").Value = 1595902.1
$ws.Range("J122").Value = 4833.8184
$ws.Range("K122").Value = 4787706.300000001
$ws.Range("L122").Value = 14501.4552
$ws.Range("M122").Value = -4785256.300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 83335816
$ws.Range("I86").Value = 2499.2
$ws.Range("J86").Value = 142859620
$ws.Range("K86").Value = 2499.2
$ws.Range("L86").Value = 142859620
$ws.Range("M86").Value = -1376.2
$ws.Range("N86").Value = -142861866

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 83335816
$ws.Range("I89").Value = 2499.2
$ws.Range("J89").Value = 142859620
$ws.Range("K89").Value = 12496
$ws.Range("L89").Value = 714298100
$ws.Range("M89").Value = -6880
$ws.Range("N89").Value = -714309332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 19444.23
$ws.Range("I99").Value = 22138.727
$ws.Range("J99").Value = 4624.5
$ws.Range("K99").Value = 22138.727
$ws.Range("L99").Value = 4624.5
$ws.Range("M99").Value = -20640.727
$ws.Range("N99").Value = -7620.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1853.4
$ws.Range("I107").Value = 1267
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1267
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 653

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14183.68
$ws.Range("I31").Value = 1502
$ws.Range("J31").Value = 33206.2
$ws.Range("K31").Value = 1502
$ws.Range("L31").Value = 33206.2
$ws.Range("M31").Value = -1207
$ws.Range("N31").Value = -33796.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 14183.68
$ws.Range("I34").Value = 1502
$ws.Range("J34").Value = 33206.2
$ws.Range("K34").Value = 1502
$ws.Range("L34").Value = 33206.2
$ws.Range("M34").Value = -1300
$ws.Range("N34").Value = -33610.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9948.673000000001
$ws.Range("I58").Value = 4420.1665
$ws.Range("J58").Value = 22387.812
$ws.Range("K58").Value = 4420.1665
$ws.Range("L58").Value = 22387.812
$ws.Range("M58").Value = -4217.1665
$ws.Range("N58").Value = -22793.812

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 15000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 15000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -17246

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 15000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 15000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -56232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1418.5294
$ws.Range("I107").Value = 1418.75
$ws.Range("J107").Value = 1418
$ws.Range("K107").Value = 1418.75
$ws.Range("L107").Value = 1418
$ws.Range("M107").Value = 501.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3982.111
$ws.Range("I132").Value = 3492.7144
$ws.Range("J132").Value = 5695
$ws.Range("K132").Value = 10478.1432
$ws.Range("L132").Value = 17085
$ws.Range("M132").Value = -7948.143199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 9948.673000000001
$ws.Range("I136").Value = 4420.1665
$ws.Range("J136").Value = 22387.812
$ws.Range("K136").Value = 13260.4995
$ws.Range("L136").Value = 67163.436
$ws.Range("M136").Value = -10710.4995
$ws.Range("N136").Value = -72263.436

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1192.4667
$ws.Range("I113").Value = 1200
$ws.Range("J113").Value = 1187.4445
$ws.Range("K113").Value = 3600
$ws.Range("L113").Value = 3562.3335
$ws.Range("M113").Value = -1430
$ws.Range("N113").Value = -7902.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1457.1
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1457.1
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 4371.299999999999
$ws.Range("N131").Value = -14451.3
$ws.Range("M131").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2220.5715
$ws.Range("I140").Value = 2220.5715
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 6661.7145
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -1481.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 10725.667
$ws.Range("I24").Value = 7296.8
$ws.Range("J24").Value = 12044.462
$ws.Range("K24").Value = 7296.8
$ws.Range("L24").Value = 12044.462
$ws.Range("M24").Value = -7123.8
$ws.Range("N24").Value = -12390.462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5633192.5
$ws.Range("I102").Value = 13517002
$ws.Range("J102").Value = 1899.8572
$ws.Range("K102").Value = 13517002
$ws.Range("L102").Value = 1899.8572
$ws.Range("M102").Value = -13515380
$ws.Range("N102").Value = -5143.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1885736.5
$ws.Range("I122").Value = 2262740.5
$ws.Range("J122").Value = 716.6667
$ws.Range("K122").Value = 6788221.5
$ws.Range("L122").Value = 2150.0001
$ws.Range("M122").Value = -6785771.5
$ws.Range("N122").Value = -7050.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4885037
$ws.Range("I126").Value = 2401351.2
$ws.Range("J126").Value = 10207220
$ws.Range("K126").Value = 7204053.600000001
$ws.Range("L126").Value = 30621660
$ws.Range("M126").Value = -7201583.600000001
$ws.Range("N126").Value = -30626600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1650663.8
$ws.Range("I7").Value = 2810134
$ws.Range("J7").Value = 8080.9165
$ws.Range("K7").Value = 2810134
$ws.Range("L7").Value = 8080.9165
$ws.Range("M7").Value = -2810022

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 38332.668
$ws.Range("I45").Value = 32499
$ws.Range("J45").Value = 50000
$ws.Range("K45").Value = 32499
$ws.Range("L45").Value = 50000
$ws.Range("M45").Value = -32092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1650663.8
$ws.Range("I126").Value = 2810134
$ws.Range("J126").Value = 8080.9165
$ws.Range("K126").Value = 8430402
$ws.Range("L126").Value = 24242.7495
$ws.Range("M126").Value = -8427932

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 960878
$ws.Range("I132").Value = 3278.2222
$ws.Range("J132").Value = 2684557.5
$ws.Range("K132").Value = 9834.6666
$ws.Range("L132").Value = 8053672.5
$ws.Range("M132").Value = -7304.6666
$ws.Range("N132").Value = -8058732.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6672354
$ws.Range("I126").Value = 5793.231
$ws.Range("J126").Value = 50005000
$ws.Range("K126").Value = 17379.693
$ws.Range("L126").Value = 150015000
$ws.Range("M126").Value = -14909.693
$ws.Range("N126").Value = -150019940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11597.75
$ws.Range("I132").Value = 3621.6667
$ws.Range("J132").Value = 99334.664
$ws.Range("K132").Value = 10865.0001
$ws.Range("L132").Value = 298003.992
$ws.Range("M132").Value = -8335.000100000001
$ws.Range("N132").Value = -303063.992
